$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill in row 6: Date, Tijd (time range), Locatie
# Reuse the existing date cell's style (B5) so no new number format is introduced
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B6").Value = 43606

$ws.Range("C6").Value = "15:30 tot 18:00"
$ws.Range("D6").Value = "E2.23"

# Update the active cell selection to E6, matching the saved view state
$ws.Range("E6").Select()

$wb.Save()
